$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add column G, copying the header style/format from column F so the
#     new header cell matches the existing bold/centered/bordered look ---
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null

# --- Header row text (column A label stays "" - header row has no A1) ---
$ws.Range("B1").Value = "Processing"
$ws.Range("C1").Value = "AI-Synonyms"
$ws.Range("D1").Value = "Product-AI"
$ws.Range("E1").Value = "Business-Process-AI"
$ws.Range("F1").Value = "Data"
$ws.Range("G1").Value = "Adjectives"

# --- Data rows: column A (the year label) is left untouched - only the
#     numeric counts in B:G change ---
$data = @(
    @{ Row = 2;  B = 49;  C = 22; D = 6;  E = 1; F = 6;  G = 25 },
    @{ Row = 3;  B = 57;  C = 24; D = 5;  E = 1; F = 10; G = 15 },
    @{ Row = 4;  B = 49;  C = 31; D = 5;  E = 0; F = 10; G = 23 },
    @{ Row = 5;  B = 50;  C = 17; D = 5;  E = 1; F = 13; G = 17 },
    @{ Row = 6;  B = 66;  C = 18; D = 10; E = 0; F = 63; G = 51 },
    @{ Row = 7;  B = 88;  C = 16; D = 3;  E = 2; F = 44; G = 52 },
    @{ Row = 8;  B = 83;  C = 31; D = 13; E = 2; F = 67; G = 109 },
    @{ Row = 9;  B = 76;  C = 28; D = 3;  E = 1; F = 41; G = 48 },
    @{ Row = 10; B = 81;  C = 28; D = 5;  E = 1; F = 55; G = 56 },
    @{ Row = 11; B = 81;  C = 29; D = 4;  E = 3; F = 60; G = 68 },
    @{ Row = 12; B = 105; C = 49; D = 4;  E = 3; F = 73; G = 85 }
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
}
